$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new text-bearing cells as Text so the numeric-looking
# strings ("0000", "00000", "000", "0") are preserved as text rather than
# being coerced to numbers, matching the original sheet's inlineStr cells.
$ws.Range("B8:G9").NumberFormat = "@"
$ws.Range("J8:J9").NumberFormat = "@"

# Copy the date-cell formatting (bold, bordered, centered, custom date
# number format) from the row above down onto the two new date cells.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# Row 8: 2025-07-31
$ws.Range("A8").Value = 45869
$ws.Range("B8").Value = "0000"
$ws.Range("C8").Value = "00000"
$ws.Range("D8").Value = "000"
$ws.Range("E8").Value = "RENDIMENTOS"
$ws.Range("F8").Value = "0"
$ws.Range("G8").Value = "0"
$ws.Range("H8").Value = 1077343.87
$ws.Range("I8").Value = "C"
$ws.Range("J8").Value = "'"
$ws.Range("K8").Value = "RENDIMENTOS"

# Row 9: 2025-08-31
$ws.Range("A9").Value = 45900
$ws.Range("B9").Value = "0000"
$ws.Range("C9").Value = "00000"
$ws.Range("D9").Value = "000"
$ws.Range("E9").Value = "RENDIMENTOS"
$ws.Range("F9").Value = "0"
$ws.Range("G9").Value = "0"
$ws.Range("H9").Value = 956265.4300000001
$ws.Range("I9").Value = "C"
$ws.Range("J9").Value = "'"
$ws.Range("K9").Value = "RENDIMENTOS"
